$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Enterprises density (per 1000 people)" row (row 11): Micro / SMEs / MSMEs
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "0.08"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "3.68"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.76"

# "Employment (% of total)" row (row 12): Micro / SMEs / MSMEs
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "7.78"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "38.97"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.75"

# "Enterprises (% of total)" row (row 14): Micro / SMEs / MSMEs
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2.01"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "95.72"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "97.73"
